$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Rename column headers
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Address"

# Update the active selection on the sheet
$ws.Activate()
$ws.Range("D6").Select()
